$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove obsolete trailing rows (35-53) - data collection window shortened
$ws.Rows("35:53").Delete()

# Update timestamps for remaining rows to reflect the new analysis run
$ws.Range("A2").Value = "05/05/2021 01:39:46"
$ws.Range("A3").Value = "05/05/2021 01:40:48"
$ws.Range("A4").Value = "05/05/2021 01:41:50"
$ws.Range("A5").Value = "05/05/2021 01:42:55"
$ws.Range("A6").Value = "05/05/2021 01:44:00"
$ws.Range("A7").Value = "05/05/2021 01:45:04"
$ws.Range("A8").Value = "05/05/2021 01:46:09"
$ws.Range("A9").Value = "05/05/2021 01:51:31"
$ws.Range("A10").Value = "05/05/2021 01:52:37"
$ws.Range("A11").Value = "05/05/2021 01:53:43"
$ws.Range("A12").Value = "05/05/2021 01:54:50"
$ws.Range("A13").Value = "05/05/2021 02:00:08"
$ws.Range("A14").Value = "05/05/2021 02:01:17"
$ws.Range("A15").Value = "05/05/2021 02:06:33"
$ws.Range("A16").Value = "05/05/2021 02:07:38"
$ws.Range("A17").Value = "05/05/2021 02:12:49"
$ws.Range("A18").Value = "05/05/2021 02:13:53"
$ws.Range("A19").Value = "05/05/2021 02:14:57"
$ws.Range("A20").Value = "05/05/2021 02:16:00"
$ws.Range("A21").Value = "05/05/2021 02:17:04"
$ws.Range("A22").Value = "05/05/2021 02:18:07"
$ws.Range("A23").Value = "05/05/2021 02:19:12"
$ws.Range("A24").Value = "05/05/2021 02:20:16"
$ws.Range("A25").Value = "05/05/2021 02:21:21"
$ws.Range("A26").Value = "05/05/2021 02:21:31"
$ws.Range("A27").Value = "05/05/2021 02:21:40"
$ws.Range("A28").Value = "05/05/2021 02:21:48"
$ws.Range("A29").Value = "05/05/2021 02:21:58"
$ws.Range("A30").Value = "05/05/2021 02:22:08"
$ws.Range("A31").Value = "05/05/2021 02:22:18"
$ws.Range("A32").Value = "05/05/2021 02:22:30"
$ws.Range("A33").Value = "05/05/2021 02:22:40"
$ws.Range("A34").Value = "05/05/2021 02:22:49"
